$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("sentences")

# Template cells for formatting reuse (existing style indices 13 and 15)
$fmtSrc13 = $ws.Range("A259:B259")
$fmtSrc15 = $ws.Range("A260:B260")

# Row 261
$ws.Range("A261").Value = 'ok'
$ws.Range("B261").Value = 'хорошо'

# Row 262
$ws.Range("A262").Value = 'that would be easy to do'
$ws.Range("B262").Value = 'это будет легко сделать'

# Row 263
$ws.Range("A263").Value = 'said fisherman'
$ws.Range("B263").Value = 'сказал рыбак'

# Row 264
$fmtSrc15.Copy()
$ws.Range("A264:B264").PasteSpecial(-4122)
$ws.Range("A264").Value = 'ok, that would be easy to do - said fisherman'
$ws.Range("B264").Value = 'хорошо, это будет легко сделать - сказал рыбак'

# Row 265
$fmtSrc15.Copy()
$ws.Range("A265:B265").PasteSpecial(-4122)
$ws.Range("A265").Value = 'the businessman nodded in agreement'
$ws.Range("B265").Value = 'бизнесмен кивнул в знак согласия'

# Row 266
$ws.Range("A266").Value = 'you can then'
$ws.Range("B266").Value = 'затем ты сможешь'

# Row 267
$ws.Range("A267").Value = 'you can then sell the extra fish'
$ws.Range("B267").Value = 'затем ты сможешь продавать дополнительную рыбу'

# Row 268
$ws.Range("A268").Value = 'and buy a bigger boat'
$ws.Range("B268").Value = 'и купить более большую лодку'

# Row 269
$ws.Range("A269").Value = 'continued the businessman'
$ws.Range("B269").Value = 'продолжил бизнесмен'

# Row 270
$fmtSrc15.Copy()
$ws.Range("A270:B270").PasteSpecial(-4122)
$ws.Range("A270").Value = 'you can then sell the extra fish and buy a bigger boat - continued businessman'
$ws.Range("B270").Value = 'затем ты сможешь продавать дополнительную рыбу и купить более большую лодку - продолжил бизнесмен'
$ws.Rows.Item(270).RowHeight = 28.8

# Row 271
$ws.Range("A271").Value = 'what for'
$ws.Range("B271").Value = 'зачем'

# Row 272
$ws.Range("A272").Value = 'asked the fisherman'
$ws.Range("B272").Value = 'спросил рыбак'

# Row 273
$ws.Range("A273").Value = 'asked the fisherman very politely'
$ws.Range("B273").Value = 'спросил рыбае очень вежливо'

# Row 274
$fmtSrc15.Copy()
$ws.Range("A274:B274").PasteSpecial(-4122)
$ws.Range("A274").Value = 'what for - asked the fisherman very politely'
$ws.Range("B274").Value = 'зачем - спросил рыбак очень вежливо'

# Row 275
$ws.Range("A275").Value = 'with a bigger boat'
$ws.Range("B275").Value = 'с более большой лодкой'

# Row 276
$ws.Range("A276").Value = 'you will catch even more fish'
$ws.Range("B276").Value = 'ты поймаешь еще больше рыбы'

# Row 277
$ws.Range("A277").Value = 'said the businessman'
$ws.Range("B277").Value = 'сказал бизнесмен'

# Row 278
$fmtSrc15.Copy()
$ws.Range("A278:B278").PasteSpecial(-4122)
$ws.Range("A278").Value = 'with a bigger boat, you will catch even more fish - said businessman'
$ws.Range("B278").Value = 'с более большой лодкой, ты поймаешь еще больше рыбы - сказал бизнесмен'

# Row 279
$ws.Range("A279").Value = 'soon'
$ws.Range("B279").Value = 'скоро'

# Row 280
$ws.Range("A280").Value = 'you will be able'
$ws.Range("B280").Value = 'ты сможешь'

# Row 281
$ws.Range("A281").Value = 'to buy another boat'
$ws.Range("B281").Value = 'купить другую лодку'

# Row 282
$ws.Range("A282").Value = 'hire people'
$ws.Range("B282").Value = 'нанять людей'

# Row 283
$ws.Range("A283").Value = 'and build a big business'
$ws.Range("B283").Value = 'и построить большой бизнес'

# Row 284
$fmtSrc15.Copy()
$ws.Range("A284:B284").PasteSpecial(-4122)
$ws.Range("A284").Value = 'soon you will be able to by another boat, hire people, and build a big business'
$ws.Range("B284").Value = 'скоро, ты сможешь купить другую лодку, нанять людей, и построить большой бизнес'
$ws.Rows.Item(284).RowHeight = 28.8

# Row 285
$fmtSrc15.Copy()
$ws.Range("A285:B285").PasteSpecial(-4122)
$ws.Range("A285").Value = 'the businessman was very excided'
$ws.Range("B285").Value = 'бизнесмен был очень взволнован'

# Row 286
$ws.Range("A286").Value = 'ones'
$ws.Range("B286").Value = 'когда'

# Row 287
$ws.Range("A287").Value = 'your business is big enough'
$ws.Range("B287").Value = 'твой бизнес станет достаточно большим'

# Row 288
$ws.Range("A288").Value = 'you can sell it'
$ws.Range("B288").Value = 'ты можешь продать его'

# Row 289
$ws.Range("A289").Value = 'and make a lot of money'
$ws.Range("B289").Value = 'и заработать много денег'

# Row 290
$fmtSrc15.Copy()
$ws.Range("A290:B290").PasteSpecial(-4122)
$ws.Range("A290").Value = 'ones your business is big enough, you can sell it and make a lot of money'
$ws.Range("B290").Value = 'когда твой бизнес станет достаточно большим, ты сможешь продать его, и заработать много денег'
$ws.Rows.Item(290).RowHeight = 28.8

# Row 291
$ws.Range("A291").Value = 'that sounds great'
$ws.Range("B291").Value = 'звучит здорово'

# Row 292
$ws.Range("A292").Value = 'said the fisherman'
$ws.Range("B292").Value = 'сказал рыбак'

# Row 293
$fmtSrc15.Copy()
$ws.Range("A293:B293").PasteSpecial(-4122)
$ws.Range("A293").Value = 'that sounds great - said the fisherman'
$ws.Range("B293").Value = 'звучит здорово - сказал рыбак'

# Row 294
$fmtSrc15.Copy()
$ws.Range("A294:B294").PasteSpecial(-4122)
$ws.Range("A294").Value = 'and what then'
$ws.Range("B294").Value = 'и что потом'

# Row 295
$ws.Range("A295").Value = 'then'
$ws.Range("B295").Value = 'затем'

# Row 296
$ws.Range("A296").Value = 'you can retire'
$ws.Range("B296").Value = 'ты сможешь уйти на пенсию'

# Row 297
$ws.Range("A297").Value = 'spend time with your family'
$ws.Range("B297").Value = 'проводить время со своей семьей'

# Row 298
$ws.Range("A298").Value = 'play the guitar'
$ws.Range("B298").Value = 'играть на гитаре'

# Row 299
$ws.Range("A299").Value = 'and enjoy life with your friends'
$ws.Range("B299").Value = 'и радоваться жизни со своими друзьями'

# Row 300
$ws.Range("A300").Value = 'said the businessman'
$ws.Range("B300").Value = 'сказал бизнесмен'

# Row 301
$fmtSrc15.Copy()
$ws.Range("A301:B301").PasteSpecial(-4122)
$ws.Range("A301").Value = 'then, you can retire, spend time with your family, play the guitar, and enjoy life with your friends - said the businessman'
$ws.Range("B301").Value = 'затем, ты можешь уйти на пенсию, проводить время с семьей, играть на гитаре, радоваться жизни со своими друзьями - сказал бизнесмен'
$ws.Rows.Item(301).RowHeight = 28.8

# Row 302
$ws.Range("A302").Value = 'the fisherman smiled and said'
$ws.Range("B302").Value = 'рыбал улыбнулся и сказал'

# Row 303
$ws.Range("A303").Value = 'is not that what I am doing right now'
$ws.Range("B303").Value = 'разве это не то, что я делаю прямо сейчас'

# Row 304
$fmtSrc15.Copy()
$ws.Range("A304:B304").PasteSpecial(-4122)
$ws.Range("A304").Value = 'the fisherman smiled and said - is not thar what I am doing right now'
$ws.Range("B304").Value = 'рыбак улыбнулся и сказал - разве это не то, что делаю прямо сейчас'

$excel.CutCopyMode = $false

# Update sheet view: scroll position + selection to match the edited range
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 284
$ws.Range("A305").Select()